$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text block (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$oldText = $wsHoja1.Range("A1").Value2
$newText = $oldText.Replace(
    "✅ 1000 Bs = 1.45 = 5464.65 pesos`n✅ 5464.65 pesos = 1.44 = 783.13 Bs",
    "✅ 1000 Bs = 1.65 = 6238.88 pesos`n✅ 6238.88 pesos = 1.64 = 921.06 Bs"
)
$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update the N10/O10/N12/O12 numeric values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 607
$wsTasas.Range("O10").Value = 3787
$wsTasas.Range("N12").Value = 3799.97
$wsTasas.Range("O12").Value = 561.001
